$wb = $excel.ActiveWorkbook

# The "Handed back" status text is shared (as a single shared string) across
# all the cells below in the original workbook, so all of them need updating
# together to keep them in sync.
$newStatus = "Handed back: not in sync with en-US"

# 1. Overview sheet: update handback status (both rows reference the same text)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# 2. zh-cn sheet: update Status column + Correspond Handback DateTime for the
#    303ebb0e file (row 2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-11-09 06:48:39"

# 3. de-de sheet: update Status column + Correspond Handback DateTime for the
#    303ebb0e file (row 2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-11-09 06:48:57"

# The longer status text no longer fits the old column width, so the
# "Status"/"Correspond Handback/Handoff File" columns widen to fit it.
$newColumnWidth = 32.62688700358076
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
